$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "System test" row of data (row 10)
$ws.Range("A10").Value = "System test"
$ws.Range("B10").Value = "Integrate the entire subsystems and test it"
$ws.Range("C10").Value = "N/A"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "N/A"
$ws.Range("F10").Value = "N/A"

# Widen column B to fit the new (longer) text
$ws.Columns.Item(2).ColumnWidth = 35.666666666666664

# Update the view: scroll so row 7 is at the top, and select D8:D9
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D8:D9").Select()
